# cau hinh luan chuyen
$wb = $excel.ActiveWorkbook

# Helper: write a value that must stay TEXT even when it looks like a
# number (Excel's COM .Value setter auto-coerces numeric-looking strings
# to real numbers). Prefixing with an apostrophe forces text entry; then
# reset the style back to "Normal" so the quote-prefix formatting Excel
# applies doesn't leave a stray style index on the cell.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Helper: make a cell a blank TEXT cell (as opposed to ClearContents,
# which leaves a typeless/empty cell). A lone leading apostrophe enters
# an empty text value; resetting the style drops the quote-prefix flag.
function Set-BlankText($range) {
    $range.Value = "'"
    $range.Style = "Normal"
}

# ---------- Sheet 1: QuyTrinh ----------
$ws1 = $wb.Worksheets.Item("QuyTrinh")

$ws1.Range("B1").Value = "Quy trình xử lý hồ sơ 1 cửa"
$ws1.Range("B2").Value = "QT 1C"
$ws1.Range("B3").Value = "H13-12345-1"

$qt1Data = @(
    @("H13-12345-1", "Quy trình xử lý hồ sơ 1 cửa", "QT 1C", "1", "Tiếp nhận hồ sơ ", "Thêm mới", "0", "Một cửa"),
    @("H13-12345-1", "Quy trình xử lý hồ sơ 1 cửa", "QT 1C", "2", "Thẩm tra tại đơn vị", "Chuyển xử lý", "10", "Chuyên viên"),
    @("H13-12345-1", "Quy trình xử lý hồ sơ 1 cửa", "QT 1C", "3", "Trình lãnh đạo phê duyệt dự thảo", "Trình phê duyệt", "5", "Lãnh đạo đơn vị"),
    @("H13-12345-1", "Quy trình xử lý hồ sơ 1 cửa", "QT 1C", "4", "Đóng dấu văn bản", "Chuyển ban hành", "5", "Lãnh đạo phòng"),
    @("H13-12345-1", "Quy trình xử lý hồ sơ 1 cửa", "QT 1C", "5", "Chuyển trả kết quả", "Chuyển trả kết quả", "0", "Một cửa")
)

$r = 6
foreach ($row in $qt1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $ws1.Cells.Item($r, 4) $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    Set-TextValue $ws1.Cells.Item($r, 7) $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------- Sheet 2: LuanChuyen ----------
$ws2 = $wb.Worksheets.Item("LuanChuyen")

$ws2.Range("B1").Value = "Quy trình xử lý hồ sơ 1 cửa"
$ws2.Range("B2").Value = "QT 1C"
$ws2.Range("B3").Value = "H13-12345-1"

$ws2.Range("A6").Value = "Tiếp nhận hồ sơ "
$ws2.Range("B6").Value = "Thẩm tra tại đơn vị"
Set-BlankText $ws2.Range("C6")
Set-BlankText $ws2.Range("D6")

$lcData = @(
    @("Thẩm tra tại đơn vị", "Trình lãnh đạo phê duyệt dự thảo", "Chuyển trả kết quả"),
    @("Trình lãnh đạo phê duyệt dự thảo", "Đóng dấu văn bản", "Thẩm tra tại đơn vị"),
    @("Đóng dấu văn bản", "Chuyển trả kết quả", "Trình lãnh đạo phê duyệt dự thảo")
)

$r = 7
foreach ($row in $lcData) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    Set-BlankText $ws2.Cells.Item($r, 4)
    $r++
}
